$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
# B2 keeps its original text type (it held the numeric-looking string "2");
# force text storage then strip the style Excel would otherwise stamp on it.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").ClearFormats()
$ws.Range("D2").Value = 0.0636
$ws.Range("G2").Value = 0.02875
$ws.Range("I2").Value = -1.36875
$ws.Range("J2").Value = -1.36875
$ws.Range("K2").Value = -2.47
$ws.Range("L2").Value = -1.54375
$ws.Range("U2").Value = 0.395
$ws.Range("V2").Value = 0.02244318181818182
$ws.Range("W2").Value = -3.416320885200554
$ws.Range("X2").Value = 0.09452487090477037
$ws.Range("Y2").Value = -3.510845756105324
$ws.Range("Z2").Value = 2.43161094224924
$ws.Range("AA2").Value = -3.328267477203648
$ws.Range("AB2").Value = 0.09452487090477037
$ws.Range("AC2").Value = -3.422792348108418
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -0.395
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.02295844231328102
$ws.Range("AK2").Value = -0.1327731092436975
$ws.Range("AL2").Value = 0.002
$ws.Range("AM2").Value = 0.002
$ws.Range("AN2").ClearContents()
$ws.Range("AO2").Value = -1095
$ws.Range("AP2").ClearContents()
$ws.Range("AQ2").Value = -1095

# Row 3 updates
$ws.Range("B3").Value = "TraceSafe Inc. (CNSX:TSF)"
$ws.Range("D3").Value = 0.0636
$ws.Range("G3").Value = 0.02875
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = -1.36875
$ws.Range("J3").Value = -1.36875
$ws.Range("K3").Value = -2.47
$ws.Range("L3").Value = -1.54375
$ws.Range("U3").Value = 0.395
$ws.Range("V3").Value = 0.02244318181818182
$ws.Range("W3").Value = -3.416320885200554
$ws.Range("X3").Value = 0.09452487090477037
$ws.Range("Y3").Value = -3.510845756105324
$ws.Range("Z3").Value = 2.43161094224924
$ws.Range("AA3").Value = -3.328267477203648
$ws.Range("AB3").Value = 0.09452487090477037
$ws.Range("AC3").Value = -3.422792348108418
$ws.Range("AG3").Value = -0.395
$ws.Range("AJ3").Value = -0.02295844231328102
$ws.Range("AK3").Value = -0.1327731092436975
$ws.Range("AL3").Value = 0.002
$ws.Range("AM3").Value = 0.002
$ws.Range("AO3").Value = -1095
$ws.Range("AQ3").Value = -1095

# Delete row 4 entirely
$ws.Rows.Item(4).Delete()
